# "long time since commit" -- re-saving the BOM sheet after a while: widen
# column A (part descriptions) and leave the view scrolled/selected where the
# author last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A needs to go from the old narrow width (~18.16 chars) to a much
# wider 74.5 (stored width units). Excel's ColumnWidth property is expressed
# in "number of characters of the Normal style font" and gets re-derived
# into the raw stored width on save, so we dial in the ColumnWidth value
# that round-trips to the target stored width of 74.5.
$ws.Columns.Item(1).ColumnWidth = 73.67

# Move the selection to where the author left off editing.
$ws.Range("C26").Select()
